function Hex2Rgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$d = $p.Designs.Item(1)
$m = $d.SlideMaster
$tcs = $m.Theme.ThemeColorScheme

$office = @("000000","FFFFFF","44546A","E7E6E6","5B9BD5","ED7D31","A5A5A5","FFC000","4472C4","70AD47","0563C1","954F72")
for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = Hex2Rgb $office[$i-1]
}
